# Insert a new survey-metadata row for Eurobarometer 98.1 (ZA7952) ahead of
# the existing 97.5 row (row 5), shifting all subsequent rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(5).Insert() | Out-Null

$ws.Range("A5").Value = "ZA7952"
# Quote-prefix so "98.1" is stored as text (matches the existing wave column
# cells, which all carry the quotePrefix style rather than being numbers).
$ws.Range("B5").Value = "'98.1"
$ws.Range("D5").Value = "Parlemeter 2022, International communications within the EU, and Key Challenges of our Times - Autumn 2022"
$ws.Range("C5").Value = "October-November 2022"

$ws.Range("C5").Select() | Out-Null
